# Scheduled runner update: refresh cached Market Board price snapshots
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ)
# and the dependent Leve price/profit columns across the job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 158.5
$ws.Range("I9").Value = 158.5
$ws.Range("K9").Value = 158.5
$ws.Range("M9").Value = 10.5

$ws.Range("H33").Value = 194.625
$ws.Range("I33").Value = 143.38461
$ws.Range("K33").Value = 143.38461
$ws.Range("M33").Value = 85.61538999999999

$ws.Range("H113").Value = 1996
$ws.Range("I113").Value = 1995
$ws.Range("K113").Value = 1995
$ws.Range("M113").Value = 1259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 509.63635
$ws.Range("I74").Value = 509.63635
$ws.Range("K74").Value = 509.63635
$ws.Range("M74").Value = 364.36365

$ws.Range("H77").Value = 509.63635
$ws.Range("I77").Value = 509.63635
$ws.Range("K77").Value = 2548.18175
$ws.Range("M77").Value = 1819.81825

$ws.Range("H102").Value = 2626
$ws.Range("I102").Value = 2335.75
$ws.Range("K102").Value = 2335.75
$ws.Range("M102").Value = -713.75

$ws.Range("H132").Value = 1937.75
$ws.Range("I132").Value = 1786
$ws.Range("K132").Value = 5358
$ws.Range("M132").Value = -2828

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2531.5417
$ws.Range("J58").Value = 5424
$ws.Range("L58").Value = 5424
$ws.Range("N58").Value = -5830

$ws.Range("H88").Value = 42183.668
$ws.Range("I88").Value = 1311
$ws.Range("J88").Value = 50358.2
$ws.Range("K88").Value = 1311
$ws.Range("L88").Value = 50358.2
$ws.Range("M88").Value = -905
$ws.Range("N88").Value = -51170.2

$ws.Range("H91").Value = 42183.668
$ws.Range("I91").Value = 1311
$ws.Range("J91").Value = 50358.2
$ws.Range("K91").Value = 1311
$ws.Range("L91").Value = 50358.2
$ws.Range("M91").Value = 93
$ws.Range("N91").Value = -53166.2

$ws.Range("H99").Value = 10990.786
$ws.Range("I99").Value = 6692.467
$ws.Range("K99").Value = 6692.467
$ws.Range("M99").Value = -5194.467

$ws.Range("H107").Value = 718.4
$ws.Range("I107").Value = 365.66666
$ws.Range("J107").Value = 1247.5
$ws.Range("K107").Value = 365.66666
$ws.Range("L107").Value = 1247.5
$ws.Range("M107").Value = 1554.33334
$ws.Range("N107").Value = -5087.5

$ws.Range("H126").Value = 10990.786
$ws.Range("I126").Value = 6692.467
$ws.Range("K126").Value = 20077.401
$ws.Range("M126").Value = -17607.401

$ws.Range("H132").Value = 2430.6667
$ws.Range("I132").Value = 2146.9
$ws.Range("K132").Value = 6440.700000000001
$ws.Range("M132").Value = -3910.700000000001

$ws.Range("H136").Value = 2531.5417
$ws.Range("J136").Value = 5424
$ws.Range("L136").Value = 16272
$ws.Range("N136").Value = -21372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 752.5714
$ws.Range("I18").Value = 711.3333
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 2133.9999
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -1964.9999
$ws.Range("N18").Value = -3338

$ws.Range("H24").Value = 1862.0834
$ws.Range("I24").Value = 237
$ws.Range("J24").Value = 2674.625
$ws.Range("K24").Value = 711
$ws.Range("L24").Value = 8023.875
$ws.Range("M24").Value = -481
$ws.Range("N24").Value = -8483.875

$ws.Range("H114").Value = 204
$ws.Range("J114").Value = 204
$ws.Range("L114").Value = 612
$ws.Range("N114").Value = -7120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4192.4443
$ws.Range("I132").Value = 3841.5
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 11524.5
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -8994.5
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6666
$ws.Range("J22").Value = 6666
$ws.Range("L22").Value = 6666
$ws.Range("N22").Value = -7256

$ws.Range("H27").Value = 6666
$ws.Range("J27").Value = 6666
$ws.Range("L27").Value = 6666
$ws.Range("N27").Value = -6880

$ws.Range("H55").Value = 538.1177
$ws.Range("I55").Value = 225.63637
$ws.Range("K55").Value = 225.63637
$ws.Range("M55").Value = -52.63637

$ws.Range("H94").Value = 25265
$ws.Range("J94").Value = 25265
$ws.Range("L94").Value = 25265
$ws.Range("N94").Value = -26617

$ws.Range("H132").Value = 2915.1667
$ws.Range("I132").Value = 2498.6667
$ws.Range("K132").Value = 7496.000100000001
$ws.Range("M132").Value = -4966.000100000001

$ws.Range("H136").Value = 8441.5
$ws.Range("I136").Value = 7995
$ws.Range("K136").Value = 23985
$ws.Range("M136").Value = -21435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20000
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20780

$ws.Range("H62").Value = 4764.684
$ws.Range("I62").Value = 3291
$ws.Range("J62").Value = 5624.3335
$ws.Range("K62").Value = 3291
$ws.Range("L62").Value = 5624.3335
$ws.Range("M62").Value = -2667
$ws.Range("N62").Value = -6872.3335

$ws.Range("H65").Value = 4764.684
$ws.Range("I65").Value = 3291
$ws.Range("J65").Value = 5624.3335
$ws.Range("K65").Value = 16455
$ws.Range("L65").Value = 28121.6675
$ws.Range("M65").Value = -13335
$ws.Range("N65").Value = -34361.6675

$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61622

$ws.Range("H69").Value = 9754
$ws.Range("J69").Value = 9754
$ws.Range("L69").Value = 9754
$ws.Range("N69").Value = -11252

$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -188112

$ws.Range("H72").Value = 9754
$ws.Range("J72").Value = 9754
$ws.Range("L72").Value = 29262
$ws.Range("N72").Value = -36750

$ws.Range("H107").Value = 180
$ws.Range("I107").Value = 180
$ws.Range("K107").Value = 540
$ws.Range("M107").Value = 1380

$ws.Range("H113").Value = 768.4
$ws.Range("I113").Value = 724.75
$ws.Range("J113").Value = 797.5
$ws.Range("K113").Value = 2174.25
$ws.Range("L113").Value = 2392.5
$ws.Range("M113").Value = -4.25
$ws.Range("N113").Value = -6732.5

$ws.Range("H126").Value = 2184.9285
$ws.Range("I126").Value = 1916.1111
$ws.Range("K126").Value = 5748.3333
$ws.Range("M126").Value = -3278.3333

$ws.Range("H132").Value = 2934.95
$ws.Range("I132").Value = 1961.7693
$ws.Range("J132").Value = 4742.2856
$ws.Range("K132").Value = 5885.3079
$ws.Range("L132").Value = 14226.8568
$ws.Range("M132").Value = -3355.3079
$ws.Range("N132").Value = -19286.8568
